$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- H8 / H9: mark as "Completado" using the same formatting already used
#     by the other "Completado" cells in column H (e.g. H10), copying only
#     the formatting (not the value) so the shared style (fill/font) is reused
#     instead of creating a brand-new cell style.
$fmtSource = $ws.Range("H10")
$fmtSource.Copy()

$ws.Range("H8").PasteSpecial(-4122)
$ws.Range("H8").Value2 = "Completado"

$ws.Range("H9").PasteSpecial(-4122)
$ws.Range("H9").Value2 = "Completado"

$excel.CutCopyMode = 0

# --- Update the "Fecha" column (I) for the affected rows from 11-Nov-2020
#     (44146) to 12-Nov-2020 (44147).
$ws.Range("I8").Value2 = 44147
$ws.Range("I9").Value2 = 44147
$ws.Range("I11").Value2 = 44147
$ws.Range("I12").Value2 = 44147
$ws.Range("I13").Value2 = 44147
$ws.Range("I14").Value2 = 44147
$ws.Range("I15").Value2 = 44147
$ws.Range("I35").Value2 = 44147
$ws.Range("I36").Value2 = 44147
$ws.Range("I37").Value2 = 44147
$ws.Range("I38").Value2 = 44147
$ws.Range("I39").Value2 = 44147

# --- Update the sheet view: scroll so the frozen/top-left cell resets to the
#     sheet default and select I35 (matches the saved selection state).
$ws.Range("I35").Select()
